# Update automàtic: dades i banners [2026-02-16 22:50]
# Refreshes the per-station weather snapshot rows (extraction timestamps
# and the latest observed measurements) on the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "2026-02-16 22:48:37"
$ws.Range("H2").Value = "'94%"
$ws.Range("N2").Value = "0.0 °C 22:16 TU"
$ws.Range("E3").Value = "2026-02-16 22:48:39"
$ws.Range("G3").Value = "246 cm"
$ws.Range("L3").Value = "73.8 km/h - 258º 22:20 TU"
$ws.Range("N3").Value = "-5.2 °C 22:11 TU"
$ws.Range("O3").Value = "-1.2 °C"
$ws.Range("E4").Value = "2026-02-16 22:48:42"
$ws.Range("H4").Value = "'63%"
$ws.Range("O4").Value = "13.1 °C"
$ws.Range("E5").Value = "2026-02-16 22:48:44"
$ws.Range("N5").Value = "-5.2 °C 22:24 TU"
$ws.Range("O5").Value = "-1.0 °C"
$ws.Range("E6").Value = "2026-02-16 22:48:46"
$ws.Range("H6").Value = "'68%"
$ws.Range("O6").Value = "11.7 °C"
$ws.Range("E7").Value = "2026-02-16 22:48:49"
$ws.Range("J7").Value = "1013.2 hPa"
$ws.Range("E8").Value = "2026-02-16 22:48:51"
$ws.Range("E9").Value = "2026-02-16 22:48:54"
$ws.Range("H9").Value = "'68%"
$ws.Range("E10").Value = "2026-02-16 22:48:56"
$ws.Range("H10").Value = "'78%"
$ws.Range("O10").Value = "10.6 °C"
$ws.Range("E11").Value = "2026-02-16 22:48:59"
$ws.Range("H11").Value = "'76%"
$ws.Range("E12").Value = "2026-02-16 22:49:01"
$ws.Range("O12").Value = "11.2 °C"
$ws.Range("E13").Value = "2026-02-16 22:49:02"
$ws.Range("E14").Value = "2026-02-16 22:49:04"
$ws.Range("E15").Value = "2026-02-16 22:49:05"
$ws.Range("E16").Value = "2026-02-16 22:49:06"
$ws.Range("N16").Value = "-5.0 °C 22:29 TU"
$ws.Range("O16").Value = "-0.4 °C"
$ws.Range("E17").Value = "2026-02-16 22:49:07"
$ws.Range("E18").Value = "2026-02-16 22:49:08"
$ws.Range("H18").Value = "'76%"
$ws.Range("E19").Value = "2026-02-16 22:49:09"
$ws.Range("K19").Value = "12.9 MJ/m2"
$ws.Range("E20").Value = "2026-02-16 22:49:10"
$ws.Range("N20").Value = "-4.0 °C 22:29 TU"
$ws.Range("O20").Value = "-0.8 °C"
$ws.Range("E21").Value = "2026-02-16 22:49:11"
$ws.Range("H21").Value = "'67%"
$ws.Range("J21").Value = "1014.0 hPa"
$ws.Range("L21").Value = "78.5 km/h - 24º 22:13 TU"
$ws.Range("O21").Value = "8.9 °C"
$ws.Range("E22").Value = "2026-02-16 22:49:12"
$ws.Range("E23").Value = "2026-02-16 22:49:15"
$ws.Range("H23").Value = "'85%"
$ws.Range("N23").Value = "-5.1 °C 22:26 TU"
$ws.Range("O23").Value = "-1.0 °C"
$ws.Range("E24").Value = "2026-02-16 22:49:17"
$ws.Range("E25").Value = "2026-02-16 22:49:20"
$ws.Range("H25").Value = "'85%"
$ws.Range("E26").Value = "2026-02-16 22:49:22"
$ws.Range("E27").Value = "2026-02-16 22:49:25"
$ws.Range("N27").Value = "-0.8 °C 22:29 TU"
$ws.Range("E28").Value = "2026-02-16 22:49:27"
$ws.Range("O28").Value = "9.5 °C"
$ws.Range("E29").Value = "2026-02-16 22:49:29"
$ws.Range("L29").Value = "25.2 km/h - 350º 22:00 TU"
$ws.Range("E30").Value = "2026-02-16 22:49:32"
$ws.Range("E31").Value = "2026-02-16 22:49:35"
$ws.Range("N31").Value = "11.1 °C 22:08 TU"
$ws.Range("E32").Value = "2026-02-16 22:49:37"
$ws.Range("E33").Value = "2026-02-16 22:49:40"
$ws.Range("H33").Value = "'72%"
$ws.Range("J33").Value = "1013.9 hPa"
$ws.Range("K33").Value = "9.8 MJ/m2"
$ws.Range("O33").Value = "6.2 °C"
$ws.Range("E34").Value = "2026-02-16 22:49:42"
$ws.Range("H34").Value = "'69%"
$ws.Range("N34").Value = "-0.5 °C 22:22 TU"
$ws.Range("O34").Value = "3.3 °C"
$ws.Range("E35").Value = "2026-02-16 22:49:45"
$ws.Range("H35").Value = "'76%"
$ws.Range("I35").Value = "3.3 mm"
$ws.Range("E36").Value = "2026-02-16 22:49:47"
$ws.Range("H36").Value = "'70%"
$ws.Range("J36").Value = "1012.7 hPa"
$ws.Range("L36").Value = "69.5 km/h - 325º 22:13 TU"
$ws.Range("O36").Value = "12.3 °C"
$ws.Range("E37").Value = "2026-02-16 22:49:50"
$ws.Range("O37").Value = "6.5 °C"
$ws.Range("E38").Value = "2026-02-16 22:49:52"
$ws.Range("E39").Value = "2026-02-16 22:49:55"
$ws.Range("H39").Value = "'77%"
$ws.Range("I39").Value = "4.8 mm"
$ws.Range("N39").Value = "-5.4 °C 22:23 TU"
$ws.Range("O39").Value = "-0.1 °C"
$ws.Range("E40").Value = "2026-02-16 22:49:57"
$ws.Range("O40").Value = "7.0 °C"
$ws.Range("E41").Value = "2026-02-16 22:50:00"
$ws.Range("J41").Value = "1014.6 hPa"
$ws.Range("K41").Value = "11.1 MJ/m2"
$ws.Range("E42").Value = "2026-02-16 22:50:02"
$ws.Range("H42").Value = "'79%"
$ws.Range("O42").Value = "11.5 °C"
$ws.Range("E43").Value = "2026-02-16 22:50:05"
$ws.Range("E44").Value = "2026-02-16 22:50:07"
$ws.Range("H44").Value = "'89%"
$ws.Range("N44").Value = "-4.6 °C 22:25 TU"
$ws.Range("O44").Value = "-0.5 °C"
$ws.Range("E45").Value = "2026-02-16 22:50:10"
$ws.Range("G45").Value = "2 cm"
$ws.Range("J45").Value = "1018.0 hPa"
$ws.Range("E46").Value = "2026-02-16 22:50:12"
